$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Id
$ws.Range("A2").Value = 61680670

# C2: Valideringsstatus
$ws.Range("C2").Value = "Behöver inte valideras"

# J2, K2, L2: new empty (text-typed) cells (Enhet, Ålder-Stadium, Kön)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = ""
$ws.Range("J2").NumberFormat = "General"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = ""
$ws.Range("K2").NumberFormat = "General"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = ""
$ws.Range("L2").NumberFormat = "General"

# P2: Lokalnamn
$ws.Range("P2").Value = "Häradssveden, Ög"

# S2: Noggrannhet
$ws.Range("S2").Value = 10

# X2: Externid (new)
$ws.Range("X2").Value = "E-Nor-0210"

# Y2: Startdatum (keep as text, not an Excel date serial)
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2016-08-11"
$ws.Range("Y2").NumberFormat = "General"

# AA2: Slutdatum (keep as text, not an Excel date serial)
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2016-08-11"
$ws.Range("AA2").NumberFormat = "General"

# AC2: Publik kommentar
$ws.Range("AC2").Value = "Området har tidigare betats av får, nu igenväxande"

# AD2: Ej återfunnen (boolean)
$ws.Range("AD2").Value = $true

# AP2: Offentlig samling (removed)
$ws.Range("AP2").ClearContents()

# AR2: Samlings-nummer (removed)
$ws.Range("AR2").ClearContents()

# AW2: Rapportör
$ws.Range("AW2").Value = "Margareta Edqvist"

# AX2: Observatörer
$ws.Range("AX2").Value = "bert lindgren, Mats Blomstedt, Rolf Wahlström"

# AY2: Projektnamn
$ws.Range("AY2").Value = "Floraväkteri Sverige"

$wb.Save()
